$d = $word.ActiveDocument

# Remove the blank "ind left=720" ListParagraph plus the four bulleted
# list items ("Un plateaux devra se déplacer...", "...apparaissent en
# cours de partie.", "Des ennemis doivent apparaître...", "Certaines
# actions ou contacts doivent faire perdre de la vie...") while keeping
# the earlier blank numId=0 list paragraph and the following
# "Lorsque le personnage..." item intact.

$startMarker = "Un plateaux devra se déplacer"
$endMarker = "nombre de vies."

$findStart = $d.Content.Duplicate
$findStart.Find.Execute($startMarker, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$findEnd = $d.Content.Duplicate
$findEnd.Find.Execute($endMarker, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Locate the 0-based paragraph indices containing the found ranges.
$targetIndex = -1
$endIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -le $findStart.Start -and $p.Range.End -ge $findStart.End) {
        $targetIndex = $i
    }
    if ($p.Range.Start -le $findEnd.Start -and $p.Range.End -ge $findEnd.End) {
        $endIndex = $i
    }
    $i = $i + 1
}

# The blank paragraph right before the "Un plateaux..." item is also
# removed, so start the deletion from it (1-based Item(targetIndex)
# is the paragraph preceding the 0-based $targetIndex paragraph).
$precedingPara = $d.Paragraphs.Item($targetIndex)
$endPara = $d.Paragraphs.Item($endIndex + 1)

$deleteRange = $d.Range($precedingPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
